$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.170.36'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '1.681.14'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.519'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +2.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.52'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.79%  '
$ws.Range('E10').Value = '  +0.65%  '
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '1.917.81'
$ws.Range('D13').Value = '1.682.09'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('E14').Value = '  +1.62%  '
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = '27.154.03'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '238.99'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').Value = '0.0₃0743'
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  +2.11%  '
$ws.Range('E23').Value = '  +3.17%  '
$ws.Range('E24').Value = '  -3.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.29'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.26'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.92%  '
$ws.Range('E28').Value = '  +0.72%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').Value = '1.570.42'
$ws.Range('E32').Value = '  +5.83%  '
$ws.Range('E33').Value = '  +1.54%  '
$ws.Range('E34').Value = '  +2.45%  '
$ws.Range('E35').Value = '  +0.71%  '
$ws.Range('E36').Value = '  +3.00%  '
$ws.Range('E37').Value = '  -1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.936'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.55%  '
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('E40').Value = '  +4.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.18'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.05%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.60'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.81%  '
$ws.Range('E44').Value = '  -2.40%  '
$ws.Range('D45').Value = '1.826.87'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.787'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.76'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('E48').Value = '  +3.07%  '
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('E50').Value = '  +6.36%  '
$ws.Range('E51').Value = '  +1.66%  '
